$wb = $excel.ActiveWorkbook

# ----- Sheet 1: VENTAS POR GRUPO -----
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Widen column B (28 -> 31 raw OOXML width units; COM ColumnWidth needs -0.83 offset)
$ws1.Columns.Item(2).ColumnWidth = 30.17

# New top record: COELLO TRONCOSO JOSE GREGORIO replaces VIEJO RIVAS MAYRA ANABELLE in row 2
$ws1.Range("B2").Value = "COELLO TRONCOSO JOSE GREGORIO"

# Insert a new row 3, pushing the old row 3 ("0 de 1" summary) down to row 4
$ws1.Rows.Item(3).Insert()

# Populate the new row 3 with the previous client's data (same as old row 2 data)
$ws1.Range("A3").Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws1.Range("B3").Value = "VIEJO RIVAS MAYRA ANABELLE"
$ws1.Range("C3:R3").Value = 0

# Update the summary label in (now) row 4 from "0 de 1" to "0 de 2"
$ws1.Range("C4:R4").Value = "0 de 2"

# ----- Sheet 2: VENTA MENSUAL -----
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

# Widen column B (28 -> 31 raw OOXML width units; COM ColumnWidth needs -0.83 offset)
$ws2.Columns.Item(2).ColumnWidth = 30.17

# New top record: COELLO TRONCOSO JOSE GREGORIO replaces VIEJO RIVAS MAYRA ANABELLE in row 2
$ws2.Range("B2").Value = "COELLO TRONCOSO JOSE GREGORIO"

# Insert a new row 3, pushing the old row 3 (totals row) down to row 4
$ws2.Rows.Item(3).Insert()

# Populate the new row 3 with the previous client's data (same as old row 2 data)
$ws2.Range("A3").Value = "ILLER LOPEZ ROBERTO FERNANDO"
$ws2.Range("B3").Value = "VIEJO RIVAS MAYRA ANABELLE"
$ws2.Range("C3:G3").Value = 0
